$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '51.943.66'
Set-TextValue 'E2' '  +0.94%  '
Set-TextValue 'D3' '2.823.51'
Set-TextValue 'E3' '  +2.67%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '355.43'
Set-TextValue 'E5' '  +6.86%  '
Set-TextValue 'D6' '113.60'
Set-TextValue 'E6' '  -1.84%  '
Set-TextValue 'D7' '0.552'
Set-TextValue 'E7' '  +2.60%  '
Set-TextValue 'E8' '  +0.08%  '
Set-TextValue 'D10' '41.77'
Set-TextValue 'E10' '  +0.49%  '
Set-TextValue 'D11' '0.0853'
Set-TextValue 'E11' '  -0.61%  '
Set-TextValue 'D12' '20.03'
Set-TextValue 'E12' '  -0.77%  '
Set-TextValue 'E13' '  +1.48%  '
Set-TextValue 'E14' '  +1.69%  '
Set-TextValue 'D15' '3.246.11'
Set-TextValue 'E15' '  +2.04%  '
Set-TextValue 'D16' '2.832.51'
Set-TextValue 'E16' '  +2.96%  '
Set-TextValue 'E17' '  +1.86%  '
Set-TextValue 'D18' '51.863.63'
Set-TextValue 'E18' '  +0.81%  '
Set-TextValue 'E19' '  +8.38%  '
Set-TextValue 'E20' '  -1.84%  '
Set-TextValue 'D21' '13.64'
Set-TextValue 'E21' '  +1.72%  '
Set-TextValue 'E22' '  +2.32%  '
Set-TextValue 'D23' '270.07'
Set-TextValue 'E23' '  -2.84%  '
Set-TextValue 'D24' '69.88'
Set-TextValue 'E24' '  +0.59%  '
Set-TextValue 'D26' '26.80'
Set-TextValue 'E26' '  +0.10%  '
Set-TextValue 'E27' '  +0.16%  '
Set-TextValue 'D28' '10.31'
Set-TextValue 'E28' '  +1.47%  '
Set-TextValue 'E29' '  +1.90%  '
Set-TextValue 'E30' '  -0.49%  '
Set-TextValue 'D31' '0.0460'
Set-TextValue 'E31' '  +33.93%  '
Set-TextValue 'D32' '50.93'
Set-TextValue 'E32' '  +2.21%  '
Set-TextValue 'D33' '33.84'
Set-TextValue 'E33' '  -3.21%  '
Set-TextValue 'D34' '5.84'
Set-TextValue 'E34' '  +5.57%  '
Set-TextValue 'D35' '0.0831'
Set-TextValue 'E35' '  +0.97%  '
Set-TextValue 'E36' '  -0.02%  '
Set-TextValue 'D38' '4.92'
Set-TextValue 'E38' '  -1.53%  '
Set-TextValue 'E39' '  -0.20%  '
Set-TextValue 'D40' '18.39'
Set-TextValue 'E40' '  -3.68%  '
Set-TextValue 'B41' 'Monero'
Set-TextValue 'C41' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D41' '128.75'
Set-TextValue 'E41' '  +1.54%  '
Set-TextValue 'B42' 'EnergySwap'
Set-TextValue 'C42' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D42' '23.53'
Set-TextValue 'E42' '  +2.41%  '
Set-TextValue 'B43' 'Stacks'
Set-TextValue 'C43' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D43' '2.56'
Set-TextValue 'E43' '  +4.77%  '
Set-TextValue 'E44' '  +1.41%  '
Set-TextValue 'D45' '2.30'
Set-TextValue 'E45' '  +0.31%  '
Set-TextValue 'D46' '3.36'
Set-TextValue 'E46' '  +1.03%  '
Set-TextValue 'D47' '2.079.97'
Set-TextValue 'E47' '  -0.47%  '
Set-TextValue 'E48' '  +3.99%  '
Set-TextValue 'E49' '  +3.16%  '
Set-TextValue 'D50' '0.940'
Set-TextValue 'E50' '  +8.64%  '
Set-TextValue 'D51' '60.66'
Set-TextValue 'E51' '  +1.42%  '

Write-Host "Applied all cell updates"
